# DaSSCo ARS Phase 3 work status overview - adjusting time estimates
# - Overview!G9 formula addend 84 -> 74 (re-estimated package G hours)
# - Overview!G11 formula addend 167 -> 120 (re-estimated package I hours)
# - Overview!J3 gets a new "Old est" header
# - Overview!J11 records the previous estimate addend (167) for reference

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Overview")

# New column header for the "old estimate" reference column
$ws.Range("J3").Value = "Old est"

# Update the two work-package hour estimates
$ws.Range("G9").Formula = "='Hours by package'!G27+74"
$ws.Range("G11").Formula = "='Hours by package'!I27+120"

# Record the previous estimate addend for package I (row 11) for reference
$ws.Range("J11").Value = 167

$wb.Application.Calculate()
